$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new numeric values (padding fix) in row 11 and row 16
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 5
$ws.Range("C16").Value = 11

# Update the view: scroll so column D is the leftmost visible column,
# and select cell G17
$ws.Range("G17").Select()
$excel.ActiveWindow.ScrollColumn = 4
